$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column W: header + formatting copied from column V ---
$ws.Range("W1").Value = "tompkins_discharged_from_hospital"
$ws.Range("V1").Copy()
$ws.Range("W1").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("V2:V30").Copy()
$ws.Range("W2:W30").PasteSpecial(-4122) # xlPasteFormats

# Fill W2:W30 with 0 (matches diff)
$ws.Range("W2:W30").Value = 0

# --- New row 31: copy formatting from row 30, then set values ---
$ws.Range("A30:W30").Copy()
$ws.Range("A31:W31").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A31").Value = 29
$ws.Range("B31").Value = 43920
$ws.Range("C31").Value = 782365
$ws.Range("D31").Value = 164566
$ws.Range("E31").Value = 37582
$ws.Range("F31").Value = 161807
$ws.Range("G31").Value = 5644
$ws.Range("H31").Value = 2978
$ws.Range("I31").Value = 66663
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1218
$ws.Range("L31").Value = 2755
$ws.Range("M31").Value = 73
$ws.Range("N31").Value = 77
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = 73
$ws.Range("Q31").Value = 1096
$ws.Range("R31").Value = 1427
$ws.Range("S31").Value = 0
$ws.Range("T31").Value = 0
$ws.Range("U31").Value = 0
$ws.Range("V31").Value = 1
$ws.Range("W31").Value = 1
